# Apply updated NATMI TPM-derived ligand-receptor statistics for Crlf1-Cntfr
# (Sending/target clusters: ECs=20, FAPs=21, MuSCs=22 in shared strings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.1993696666666667
$ws.Range("H2").Value = 0.598109
$ws.Range("I2").Value = 0.007341795645341
$ws.Range("J2").Value = 0.007341795645341001
$ws.Range("M2").Value = 0.04154133333333333
$ws.Range("N2").Value = 0.124624
$ws.Range("O2").Value = 0.002517093804502335
$ws.Range("P2").Value = 0.002517093804502335
$ws.Range("Q2").Value = 0.008282081779555556
$ws.Range("R2").Value = 0.074538736016
$ws.Range("S2").Value = 0.00001847998833281005
$ws.Range("T2").Value = 0.00001847998833281005

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.1993696666666667
$ws.Range("H3").Value = 0.598109
$ws.Range("I3").Value = 0.007341795645341
$ws.Range("J3").Value = 0.007341795645341001
$ws.Range("O3").Value = 0.8133013372545576
$ws.Range("P3").Value = 0.8133013372545578
$ws.Range("Q3").Value = 2.676033834939222
$ws.Range("R3").Value = 24.084304514453
$ws.Range("S3").Value = 0.005971092216205523
$ws.Range("T3").Value = 0.005971092216205525

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.1993696666666667
$ws.Range("H4").Value = 0.598109
$ws.Range("I4").Value = 0.007341795645341
$ws.Range("J4").Value = 0.007341795645341001
$ws.Range("O4").Value = 0.18418156894094
$ws.Range("P4").Value = 0.18418156894094
$ws.Range("Q4").Value = 0.6060190579815555
$ws.Range("R4").Value = 5.454171521834
$ws.Range("S4").Value = 0.001352223440802666
$ws.Range("T4").Value = 0.001352223440802667

# Row 5
$ws.Range("I5").Value = 0.6504066282123248
$ws.Range("J5").Value = 0.6504066282123249
$ws.Range("M5").Value = 0.04154133333333333
$ws.Range("N5").Value = 0.124624
$ws.Range("O5").Value = 0.002517093804502335
$ws.Range("P5").Value = 0.002517093804502335
$ws.Range("Q5").Value = 0.7337061864746666
$ws.Range("R5").Value = 6.603355678272
$ws.Range("S5").Value = 0.001637134494280496
$ws.Range("T5").Value = 0.001637134494280496

# Row 6
$ws.Range("I6").Value = 0.6504066282123248
$ws.Range("J6").Value = 0.6504066282123249
$ws.Range("O6").Value = 0.8133013372545576
$ws.Range("P6").Value = 0.8133013372545578
$ws.Range("S6").Value = 0.5289765804843115
$ws.Range("T6").Value = 0.5289765804843118

# Row 7
$ws.Range("I7").Value = 0.6504066282123248
$ws.Range("J7").Value = 0.6504066282123249
$ws.Range("O7").Value = 0.18418156894094
$ws.Range("P7").Value = 0.18418156894094
$ws.Range("S7").Value = 0.1197929132337326
$ws.Range("T7").Value = 0.1197929132337326

# Row 8
$ws.Range("I8").Value = 0.3422515761423342
$ws.Range("J8").Value = 0.3422515761423342
$ws.Range("M8").Value = 0.04154133333333333
$ws.Range("N8").Value = 0.124624
$ws.Range("O8").Value = 0.002517093804502335
$ws.Range("P8").Value = 0.002517093804502335
$ws.Range("Q8").Value = 0.386084778128
$ws.Range("R8").Value = 3.474763003152
$ws.Range("S8").Value = 0.0008614793218890285
$ws.Range("T8").Value = 0.0008614793218890285

# Row 9
$ws.Range("I9").Value = 0.3422515761423342
$ws.Range("J9").Value = 0.3422515761423342
$ws.Range("O9").Value = 0.8133013372545576
$ws.Range("P9").Value = 0.8133013372545578
$ws.Range("S9").Value = 0.2783536645540404
$ws.Range("T9").Value = 0.2783536645540405

# Row 10
$ws.Range("I10").Value = 0.3422515761423342
$ws.Range("J10").Value = 0.3422515761423342
$ws.Range("O10").Value = 0.18418156894094
$ws.Range("P10").Value = 0.18418156894094
$ws.Range("S10").Value = 0.0630364322664047
$ws.Range("T10").Value = 0.06303643226640471

